$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 with trade data, mirroring the structure of existing rows
$ws.Cells.Item(12, 1).Value = 9843.91
$ws.Cells.Item(12, 2).Value = 9831.1299999999992
$ws.Cells.Item(12, 3).Value = 77.78
$ws.Cells.Item(12, 4).Value = 77.88
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = 0.13
$ws.Cells.Item(12, 7).Value = 42620.766134259262
$ws.Cells.Item(12, 8).Value = $true

# Apply the same number format style as the other Date column cells (style index 1 -> numFmtId 22)
$ws.Cells.Item(11, 7).Copy() | Out-Null
$ws.Cells.Item(12, 7).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12, 7).Value = 42620.766134259262
